# Generate Report for Handoff
# Updates status strings "In Translation" -> "Ready for handoff" and
# refreshes the handoff timestamps on the Overview / zh-cn / de-de sheets.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Overview sheet: B2 = zh-cn status, C2 = de-de status, D2 = Latest Handoff Date
$overview.Range("B2").Value = "Ready for handoff"
$overview.Range("C2").Value = "Ready for handoff"
$overview.Range("D2").Value = "2016-26-12 12:26:52"

# zh-cn sheet: C2 = Status, E2 = Latest Handoff Datetime
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("E2").Value = "2016-03-12 12:26:48"

# de-de sheet: C2 = Status, E2 = Latest Handoff Datetime
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("E2").Value = "2016-03-12 12:26:52"
